$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1, B2, B3 previously held the shared-string value "이김" (win). Replace
# them with the literal numeric score 10000 (the intended score update
# referenced in the commit message).
$ws.Range("B1").Value = 10000
$ws.Range("B2").Value = 10000
$ws.Range("B3").Value = 10000
